# Update countries & provincias Spain
# - Reorder "Pakistan"/"Malasia" and "Isla de Man"/"Guadalupe" (rows swap position)
# - Refresh case numbers for several countries
# - Bump the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Pakistan / Malasia (row 33 <-> row 34) ---
# Row 33 becomes Pakistan with refreshed figures
$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 4004
$ws.Range("C33").Value = 238
$ws.Range("D33").Value = 429
$ws.Range("E33").Value = 3521
$ws.Range("F33").Value = 28
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 54

# Row 34 becomes Malasia, carrying the old (unchanged) Malasia figures
$ws.Range("A34").Value = "Malasia"
$ws.Range("B34").Value = 3963
$ws.Range("C34").Value = 170
$ws.Range("D34").Value = 1321
$ws.Range("E34").Value = 2579
$ws.Range("F34").Value = 92
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 63

# --- Swap Isla de Man / Guadalupe (row 120 <-> row 121) ---
# Row 120 becomes Isla de Man with refreshed figures
$ws.Range("A120").Value = "Isla de Man"
$ws.Range("B120").Value = 150
$ws.Range("C120").Value = 11
$ws.Range("D120").Value = 73
$ws.Range("E120").Value = 76
$ws.Range("F120").Value = 6
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1

# Row 121 becomes Guadalupe, carrying the old (unchanged) Guadalupe figures
$ws.Range("A121").Value = "Guadalupe"
$ws.Range("B121").Value = 139
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 31
$ws.Range("E121").Value = 101
$ws.Range("F121").Value = 14
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 7

# --- Refresh case numbers for other countries ---

# Row 10: Iran
$ws.Range("B10").Value = 62589
$ws.Range("C10").Value = 2089
$ws.Range("D10").Value = 27039
$ws.Range("E10").Value = 31678
$ws.Range("F10").Value = 3987
$ws.Range("G10").Value = 133
$ws.Range("H10").Value = 3872

# Row 17: Austria
$ws.Range("B17").Value = 12427
$ws.Range("C17").Value = 130
$ws.Range("E17").Value = 8138

# Row 32: Rumania
$ws.Range("B32").Value = 4417
$ws.Range("C32").Value = 360
$ws.Range("D32").Value = 460
$ws.Range("E32").Value = 3775

# Row 75: Kazajistan
$ws.Range("D75").Value = 50
$ws.Range("E75").Value = 629

# Row 90: Albania
$ws.Range("E90").Value = 230
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 22

# Row 104: Vietnam
$ws.Range("D104").Value = 123
$ws.Range("E104").Value = 122

# --- Bump the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 12:22"
